# Raport.xlsx edit: record two new finished test files (LekarzControllerTest.cs,
# OsobaControllerTest.cs) in the next free "Data / Plik / Linie" column group
# (N/O/P) for the existing rows 23 and 24, matching the pattern already used
# on row 22 (and rows 7-22 above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$loggedDate = (Get-Date -Year 2025 -Month 5 -Day 16 -Hour 0 -Minute 0 -Second 0).Date

# Row 23: LekarzControllerTest.cs, 90 lines, logged 2025-05-16
$ws.Range("N23").Value = $loggedDate
$ws.Range("O23").Value = "LekarzControllerTest.cs"
$ws.Range("P23").Value = 90

# Row 24: OsobaControllerTest.cs, 126 lines, logged 2025-05-16
$ws.Range("N24").Value = $loggedDate
$ws.Range("O24").Value = "OsobaControllerTest.cs"
$ws.Range("P24").Value = 126

# Move the active selection to where the user's cursor ended up after the edit.
$ws.Range("O31").Select()
